$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.023.85"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.299.27"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.99%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.99%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "2.656.81"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "2.287.31"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").Value = "42.924.82"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.86%  "
$ws.Range("D20").Value = "0.0₃0903"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.93%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "169.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").Value = "  -10.31%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0689"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("D44").Value = "1.985.54"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("E49").Value = "  +3.74%  "
$ws.Range("D50").Value = "2.522.99"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("E51").Value = "  +0.47%  "
